$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Powder Keg entry. Typed order: A12, C12, B12
$ws.Range("A12").Value = "ATA_ITEM_POWDER_KEG"
$ws.Range("C12").Value = "火药桶"
$ws.Range("B12").Value = "Powder Keg"

# Rows 13-18: Rum entries. Typed order: B13, C13, then A13..A18 top to bottom
$ws.Range("B13").Value = "Rum"
$ws.Range("C13").Value = "朗姆酒"

$ws.Range("A13").Value = "ATA_ITEM_RUM_5"
$ws.Range("A14").Value = "ATA_ITEM_RUM_4"
$ws.Range("A15").Value = "ATA_ITEM_RUM_3"
$ws.Range("A16").Value = "ATA_ITEM_RUM_2"
$ws.Range("A17").Value = "ATA_ITEM_RUM_1"
$ws.Range("A18").Value = "ATA_ITEM_RUM_0"

# Fill remaining B/C cells for rows 14-18 (reuse existing shared strings)
$ws.Range("B14").Value = "Rum"
$ws.Range("C14").Value = "朗姆酒"
$ws.Range("B15").Value = "Rum"
$ws.Range("C15").Value = "朗姆酒"
$ws.Range("B16").Value = "Rum"
$ws.Range("C16").Value = "朗姆酒"
$ws.Range("B17").Value = "Rum"
$ws.Range("C17").Value = "朗姆酒"
$ws.Range("B18").Value = "Rum"
$ws.Range("C18").Value = "朗姆酒"

[void]$ws.Range("C13").Select()
